$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values (formerly "Strike#"); regenerate with new computed values
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
